$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.26
$ws.Range("F2").Value = 0.89

$ws.Range("B3").Value = 1.58
$ws.Range("F3").Value = 1.23

$ws.Range("C4").Value = 1.46

$ws.Range("C5").Value = 1.39
$ws.Range("D5").Value = 1.34
$ws.Range("G5").Value = 0.74

$ws.Range("B6").Value = 1.95
$ws.Range("C6").Value = 1.47

$ws.Range("G7").Value = 1.16
